$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# B2: Date - use a leading apostrophe so Excel keeps it as literal text
# instead of auto-converting the "yyyy-mm-dd" string into a date serial.
$ws.Range("B2").Value = "'2025-07-11"

# B5: Inclusion Criteria bullet list
$inclusionCriteria = "`n    " + [char]0x2022 + " Studies published in English, peer-reviewed journals" + `
    "`n    " + [char]0x2022 + " About leptin and Alzheimer" + [char]0x2019 + "s" + `
    "`n    " + [char]0x2022 + " Relevant papers available as full text" + `
    "`n    " + [char]0x2022 + " Randomized control trials " + `
    "`n    "
$ws.Range("B5").Value = $inclusionCriteria

# B6: Stopping Criteria
$ws.Range("B6").Value = "40% of total quota selected for tranche"

# B7: Study Type
$ws.Range("B7").Value = "Randomized control trials"
